$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$ws.Cells.Item(50, 1).Value = '2024-09-29 03:09:18'
$ws.Cells.Item(50, 2).Value = 'monitor_price'
$ws.Cells.Item(50, 3).Value = 'https://example.com/product'
Set-TextCell 50 4 '$199.99'
Set-TextCell 50 5 '2024-09-29'
$ws.Cells.Item(50, 6).Value = '03:09:18'

$ws.Cells.Item(51, 1).Value = '2024-09-29 03:09:18'
$ws.Cells.Item(51, 2).Value = 'monitor_price'
$ws.Cells.Item(51, 3).Value = 'invalid_url'
Set-TextCell 51 4 'Error fetching price: Invalid URL'
Set-TextCell 51 5 '2024-09-29'
$ws.Cells.Item(51, 6).Value = '03:09:18'

$ws.Cells.Item(52, 1).Value = '2024-09-29 03:15:00'
$ws.Cells.Item(52, 2).Value = 'monitor_price'
$ws.Cells.Item(52, 3).Value = 'https://example.com/product'
Set-TextCell 52 4 '$199.99'
Set-TextCell 52 5 '2024-09-29'
$ws.Cells.Item(52, 6).Value = '03:15:00'

$ws.Cells.Item(53, 1).Value = '2024-09-29 03:15:01'
$ws.Cells.Item(53, 2).Value = 'monitor_price'
$ws.Cells.Item(53, 3).Value = 'https://example.com/product'
Set-TextCell 53 4 '$199.99'
Set-TextCell 53 5 '2024-09-29'
$ws.Cells.Item(53, 6).Value = '03:15:01'

$ws.Cells.Item(54, 1).Value = '2024-09-29 03:16:05'
$ws.Cells.Item(54, 2).Value = 'monitor_price'
$ws.Cells.Item(54, 3).Value = 'https://example.com/product'
Set-TextCell 54 4 '$199.99'
Set-TextCell 54 5 '2024-09-29'
$ws.Cells.Item(54, 6).Value = '03:16:05'

$ws.Cells.Item(55, 1).Value = '2024-09-29 03:16:06'
$ws.Cells.Item(55, 2).Value = 'monitor_price'
$ws.Cells.Item(55, 3).Value = 'https://example.com/product'
Set-TextCell 55 4 '$199.99'
Set-TextCell 55 5 '2024-09-29'
$ws.Cells.Item(55, 6).Value = '03:16:06'

$ws.Cells.Item(56, 1).Value = '2024-09-29 03:17:25'
$ws.Cells.Item(56, 2).Value = 'monitor_price'
$ws.Cells.Item(56, 3).Value = 'https://example.com/product'
Set-TextCell 56 4 '$199.99'
Set-TextCell 56 5 '2024-09-29'
$ws.Cells.Item(56, 6).Value = '03:17:25'

$ws.Cells.Item(57, 1).Value = '2024-09-29 03:17:26'
$ws.Cells.Item(57, 2).Value = 'monitor_price'
$ws.Cells.Item(57, 3).Value = 'https://example.com/product'
Set-TextCell 57 4 '$199.99'
Set-TextCell 57 5 '2024-09-29'
$ws.Cells.Item(57, 6).Value = '03:17:26'

$ws.Cells.Item(58, 1).Value = '2024-09-29 03:21:01'
$ws.Cells.Item(58, 2).Value = 'monitor_price'
$ws.Cells.Item(58, 3).Value = 'https://example.com/product'
Set-TextCell 58 4 '$199.99'
Set-TextCell 58 5 '2024-09-29'
$ws.Cells.Item(58, 6).Value = '03:21:01'

$ws.Cells.Item(59, 1).Value = '2024-09-29 03:21:03'
$ws.Cells.Item(59, 2).Value = 'monitor_price'
$ws.Cells.Item(59, 3).Value = 'https://example.com/product'
Set-TextCell 59 4 '$199.99'
Set-TextCell 59 5 '2024-09-29'
$ws.Cells.Item(59, 6).Value = '03:21:03'

$ws.Cells.Item(60, 1).Value = '2024-09-29 03:21:04'
$ws.Cells.Item(60, 2).Value = 'monitor_price'
$ws.Cells.Item(60, 3).Value = 'https://example.com/product'
Set-TextCell 60 4 'Failed to fetch price: Error on page'
Set-TextCell 60 5 '2024-09-29'
$ws.Cells.Item(60, 6).Value = '03:21:04'

$ws.Cells.Item(61, 1).Value = '2024-09-29 03:23:31'
$ws.Cells.Item(61, 2).Value = 'monitor_price'
$ws.Cells.Item(61, 3).Value = 'https://example.com/product'
Set-TextCell 61 4 '100 USD'
Set-TextCell 61 5 '2024-09-29'
$ws.Cells.Item(61, 6).Value = '03:23:31'

$ws.Cells.Item(62, 1).Value = '2024-09-29 03:23:32'
$ws.Cells.Item(62, 2).Value = 'monitor_price'
$ws.Cells.Item(62, 3).Value = 'https://example.com/product'
Set-TextCell 62 4 '100 USD'
Set-TextCell 62 5 '2024-09-29'
$ws.Cells.Item(62, 6).Value = '03:23:32'

$ws.Cells.Item(63, 1).Value = '2024-09-29 03:23:33'
$ws.Cells.Item(63, 2).Value = 'monitor_price'
$ws.Cells.Item(63, 3).Value = 'https://example.com/product'
Set-TextCell 63 4 '100 USD'
Set-TextCell 63 5 '2024-09-29'
$ws.Cells.Item(63, 6).Value = '03:23:33'

$ws.Cells.Item(64, 1).Value = '2024-09-29 03:23:34'
$ws.Cells.Item(64, 2).Value = 'monitor_price'
$ws.Cells.Item(64, 3).Value = 'https://example.com/product'
Set-TextCell 64 4 '100 USD'
Set-TextCell 64 5 '2024-09-29'
$ws.Cells.Item(64, 6).Value = '03:23:34'

$ws.Cells.Item(65, 1).Value = '2024-09-29 03:23:35'
$ws.Cells.Item(65, 2).Value = 'monitor_price'
$ws.Cells.Item(65, 3).Value = 'https://example.com/product'
Set-TextCell 65 4 '100 USD'
Set-TextCell 65 5 '2024-09-29'
$ws.Cells.Item(65, 6).Value = '03:23:35'

$ws.Cells.Item(66, 1).Value = '2024-09-29 03:23:36'
$ws.Cells.Item(66, 2).Value = 'monitor_price'
$ws.Cells.Item(66, 3).Value = 'https://example.com/product'
Set-TextCell 66 4 '100 USD'
Set-TextCell 66 5 '2024-09-29'
$ws.Cells.Item(66, 6).Value = '03:23:36'

$ws.Cells.Item(67, 1).Value = '2024-09-29 03:23:37'
$ws.Cells.Item(67, 2).Value = 'monitor_price'
$ws.Cells.Item(67, 3).Value = 'https://example.com/product'
Set-TextCell 67 4 '100 USD'
Set-TextCell 67 5 '2024-09-29'
$ws.Cells.Item(67, 6).Value = '03:23:37'

$ws.Cells.Item(68, 1).Value = '2024-09-29 03:23:38'
$ws.Cells.Item(68, 2).Value = 'monitor_price'
$ws.Cells.Item(68, 3).Value = 'https://example.com/product'
Set-TextCell 68 4 '100 USD'
Set-TextCell 68 5 '2024-09-29'
$ws.Cells.Item(68, 6).Value = '03:23:38'

$ws.Cells.Item(69, 1).Value = '2024-09-29 03:23:39'
$ws.Cells.Item(69, 2).Value = 'monitor_price'
$ws.Cells.Item(69, 3).Value = 'https://example.com/product'
Set-TextCell 69 4 '100 USD'
Set-TextCell 69 5 '2024-09-29'
$ws.Cells.Item(69, 6).Value = '03:23:39'

$ws.Cells.Item(70, 1).Value = '2024-09-29 03:23:41'
$ws.Cells.Item(70, 2).Value = 'monitor_price'
$ws.Cells.Item(70, 3).Value = 'https://example.com/product'
Set-TextCell 70 4 '100 USD'
Set-TextCell 70 5 '2024-09-29'
$ws.Cells.Item(70, 6).Value = '03:23:41'

$ws.Cells.Item(71, 1).Value = '2024-09-29 03:23:42'
$ws.Cells.Item(71, 2).Value = 'monitor_price'
$ws.Cells.Item(71, 3).Value = 'https://example.com/product'
Set-TextCell 71 4 '100 USD'
Set-TextCell 71 5 '2024-09-29'
$ws.Cells.Item(71, 6).Value = '03:23:42'

$ws.Cells.Item(72, 1).Value = '2024-09-29 03:23:43'
$ws.Cells.Item(72, 2).Value = 'monitor_price'
$ws.Cells.Item(72, 3).Value = 'https://example.com/product'
Set-TextCell 72 4 '100 USD'
Set-TextCell 72 5 '2024-09-29'
$ws.Cells.Item(72, 6).Value = '03:23:43'

$ws.Cells.Item(73, 1).Value = '2024-09-29 03:23:44'
$ws.Cells.Item(73, 2).Value = 'monitor_price'
$ws.Cells.Item(73, 3).Value = 'https://example.com/product'
Set-TextCell 73 4 '100 USD'
Set-TextCell 73 5 '2024-09-29'
$ws.Cells.Item(73, 6).Value = '03:23:44'

$ws.Cells.Item(74, 1).Value = '2024-09-29 03:23:45'
$ws.Cells.Item(74, 2).Value = 'monitor_price'
$ws.Cells.Item(74, 3).Value = 'https://example.com/product'
Set-TextCell 74 4 '100 USD'
Set-TextCell 74 5 '2024-09-29'
$ws.Cells.Item(74, 6).Value = '03:23:45'

$ws.Cells.Item(75, 1).Value = '2024-09-29 03:23:46'
$ws.Cells.Item(75, 2).Value = 'monitor_price'
$ws.Cells.Item(75, 3).Value = 'https://example.com/product'
Set-TextCell 75 4 '100 USD'
Set-TextCell 75 5 '2024-09-29'
$ws.Cells.Item(75, 6).Value = '03:23:46'

$ws.Cells.Item(76, 1).Value = '2024-09-29 03:23:47'
$ws.Cells.Item(76, 2).Value = 'monitor_price'
$ws.Cells.Item(76, 3).Value = 'https://example.com/product'
Set-TextCell 76 4 '100 USD'
Set-TextCell 76 5 '2024-09-29'
$ws.Cells.Item(76, 6).Value = '03:23:47'

$ws.Cells.Item(77, 1).Value = '2024-09-29 03:23:48'
$ws.Cells.Item(77, 2).Value = 'monitor_price'
$ws.Cells.Item(77, 3).Value = 'https://example.com/product'
Set-TextCell 77 4 '100 USD'
Set-TextCell 77 5 '2024-09-29'
$ws.Cells.Item(77, 6).Value = '03:23:48'

$ws.Cells.Item(78, 1).Value = '2024-09-29 03:23:49'
$ws.Cells.Item(78, 2).Value = 'monitor_price'
$ws.Cells.Item(78, 3).Value = 'https://example.com/product'
Set-TextCell 78 4 '100 USD'
Set-TextCell 78 5 '2024-09-29'
$ws.Cells.Item(78, 6).Value = '03:23:49'

$ws.Cells.Item(79, 1).Value = '2024-09-29 03:23:50'
$ws.Cells.Item(79, 2).Value = 'monitor_price'
$ws.Cells.Item(79, 3).Value = 'https://example.com/product'
Set-TextCell 79 4 '100 USD'
Set-TextCell 79 5 '2024-09-29'
$ws.Cells.Item(79, 6).Value = '03:23:50'

$ws.Cells.Item(80, 1).Value = '2024-09-29 03:23:51'
$ws.Cells.Item(80, 2).Value = 'monitor_price'
$ws.Cells.Item(80, 3).Value = 'https://example.com/product'
Set-TextCell 80 4 '100 USD'
Set-TextCell 80 5 '2024-09-29'
$ws.Cells.Item(80, 6).Value = '03:23:51'

$ws.Cells.Item(81, 1).Value = '2024-09-29 03:23:52'
$ws.Cells.Item(81, 2).Value = 'monitor_price'
$ws.Cells.Item(81, 3).Value = 'https://example.com/product'
Set-TextCell 81 4 '100 USD'
Set-TextCell 81 5 '2024-09-29'
$ws.Cells.Item(81, 6).Value = '03:23:52'

$ws.Cells.Item(82, 1).Value = '2024-09-29 03:23:53'
$ws.Cells.Item(82, 2).Value = 'monitor_price'
$ws.Cells.Item(82, 3).Value = 'https://example.com/product'
Set-TextCell 82 4 '100 USD'
Set-TextCell 82 5 '2024-09-29'
$ws.Cells.Item(82, 6).Value = '03:23:53'

$ws.Cells.Item(83, 1).Value = '2024-09-29 03:23:54'
$ws.Cells.Item(83, 2).Value = 'monitor_price'
$ws.Cells.Item(83, 3).Value = 'https://example.com/product'
Set-TextCell 83 4 '100 USD'
Set-TextCell 83 5 '2024-09-29'
$ws.Cells.Item(83, 6).Value = '03:23:54'

$ws.Cells.Item(84, 1).Value = '2024-09-29 03:23:55'
$ws.Cells.Item(84, 2).Value = 'monitor_price'
$ws.Cells.Item(84, 3).Value = 'https://example.com/product'
Set-TextCell 84 4 '100 USD'
Set-TextCell 84 5 '2024-09-29'
$ws.Cells.Item(84, 6).Value = '03:23:55'

$ws.Cells.Item(85, 1).Value = '2024-09-29 03:23:56'
$ws.Cells.Item(85, 2).Value = 'monitor_price'
$ws.Cells.Item(85, 3).Value = 'https://example.com/product'
Set-TextCell 85 4 '100 USD'
Set-TextCell 85 5 '2024-09-29'
$ws.Cells.Item(85, 6).Value = '03:23:56'

$ws.Cells.Item(86, 1).Value = '2024-09-29 03:23:57'
$ws.Cells.Item(86, 2).Value = 'monitor_price'
$ws.Cells.Item(86, 3).Value = 'https://example.com/product'
Set-TextCell 86 4 '100 USD'
Set-TextCell 86 5 '2024-09-29'
$ws.Cells.Item(86, 6).Value = '03:23:57'

$ws.Cells.Item(87, 1).Value = '2024-09-29 03:23:58'
$ws.Cells.Item(87, 2).Value = 'monitor_price'
$ws.Cells.Item(87, 3).Value = 'https://example.com/product'
Set-TextCell 87 4 '100 USD'
Set-TextCell 87 5 '2024-09-29'
$ws.Cells.Item(87, 6).Value = '03:23:58'

$ws.Cells.Item(88, 1).Value = '2024-09-29 03:24:00'
$ws.Cells.Item(88, 2).Value = 'monitor_price'
$ws.Cells.Item(88, 3).Value = 'https://example.com/product'
Set-TextCell 88 4 '100 USD'
Set-TextCell 88 5 '2024-09-29'
$ws.Cells.Item(88, 6).Value = '03:24:00'

$ws.Cells.Item(89, 1).Value = '2024-09-29 03:24:01'
$ws.Cells.Item(89, 2).Value = 'monitor_price'
$ws.Cells.Item(89, 3).Value = 'https://example.com/product'
Set-TextCell 89 4 '100 USD'
Set-TextCell 89 5 '2024-09-29'
$ws.Cells.Item(89, 6).Value = '03:24:01'

$ws.Cells.Item(90, 1).Value = '2024-09-29 03:24:02'
$ws.Cells.Item(90, 2).Value = 'monitor_price'
$ws.Cells.Item(90, 3).Value = 'https://example.com/product'
Set-TextCell 90 4 '100 USD'
Set-TextCell 90 5 '2024-09-29'
$ws.Cells.Item(90, 6).Value = '03:24:02'

$ws.Cells.Item(91, 1).Value = '2024-09-29 03:25:36'
$ws.Cells.Item(91, 2).Value = 'monitor_price'
$ws.Cells.Item(91, 3).Value = 'https://example.com/product'
Set-TextCell 91 4 '100 USD'
Set-TextCell 91 5 '2024-09-29'
$ws.Cells.Item(91, 6).Value = '03:25:36'

$ws.Cells.Item(92, 1).Value = '2024-09-29 03:26:34'
$ws.Cells.Item(92, 2).Value = 'monitor_price'
$ws.Cells.Item(92, 3).Value = 'https://example.com/product'
Set-TextCell 92 4 '100 USD'
Set-TextCell 92 5 '2024-09-29'
$ws.Cells.Item(92, 6).Value = '03:26:34'
